$d = $word.ActiveDocument

# Insert four new lyric paragraphs ("Veja! ..." / blank / "Tente! ..." / blank)
# ahead of all existing content, styled in dark blue (theme accent1, shade 80%)
# at 16pt (sz/szCs = 32 half-points), matching the target OOXML exactly -
# including the bold lead-in words and the gramStart/gramEnd proofing marks
# around "tá". Using Range.InsertXML (a WordprocessingML package fragment)
# lets us set w:themeColor/w:themeShade and w:proofErr precisely, rather than
# relying on Font property coercion.
$insertionPoint = $d.Range($d.Content.Start, $d.Content.Start)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:color w:val="1F3864" w:themeColor="accent1" w:themeShade="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="1F3864" w:themeColor="accent1" w:themeShade="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Veja</w:t></w:r><w:r><w:rPr><w:color w:val="1F3864" w:themeColor="accent1" w:themeShade="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">! Não diga que a canção está perdida, tenha fé em Deus, tenha fé na vida, tente outra vez. Beba, pois a água viva ainda </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="1F3864" w:themeColor="accent1" w:themeShade="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>tá</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="1F3864" w:themeColor="accent1" w:themeShade="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> na fonte, você tem dois pés para cruzar a ponte, nada acabou.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="1F3864" w:themeColor="accent1" w:themeShade="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="1F3864" w:themeColor="accent1" w:themeShade="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="1F3864" w:themeColor="accent1" w:themeShade="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Tente</w:t></w:r><w:r><w:rPr><w:color w:val="1F3864" w:themeColor="accent1" w:themeShade="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>! Levante sua mão sedenta e recomece a andar, não pense que a cabeça aguenta se você parar. Há uma voz que canta, há uma voz que dança, uma voz que gira bailando no ar.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="1F3864" w:themeColor="accent1" w:themeShade="80"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xml)
